# Applies the cryptos.xlsx price/volume refresh described in the commit
# "Updated cryptos list ... with GitHub Actions", including the PEPE/
# PancakeSwap and Hedera/EthereumClassic row swaps (rows 31-34).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells (column D) are stored as text (e.g. '65.482.81') in this
# sheet. Plain .Value assignment lets Excel auto-convert parseable
# numeric-looking strings into real numbers, so for those cells we
# temporarily force a text number format, assign the value, then restore
# the default 'Normal' style so no stray formatting is left behind.
function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "65.482.81"
$ws.Range("E2").Value = "  +6.61%  "
$ws.Range("D3").Value = "2.995.74"
$ws.Range("E3").Value = "  +3.69%  "
$ws.Range("E4").Value = "  +0.17%  "
Set-TextValue "D5" "584.79"
$ws.Range("E5").Value = "  +3.01%  "
Set-TextValue "D6" "153.48"
$ws.Range("E6").Value = "  +6.78%  "
Set-TextValue "D7" "1.00"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +2.15%  "
$ws.Range("D9").Value = "2.991.25"
$ws.Range("E9").Value = "  +3.55%  "
Set-TextValue "D10" "6.96"
$ws.Range("E10").Value = "  -1.15%  "
Set-TextValue "D11" "0.152"
$ws.Range("E11").Value = "  +4.05%  "
$ws.Range("E12").Value = "  +3.78%  "
$ws.Range("E13").Value = "  +2.89%  "
Set-TextValue "D14" "33.91"
$ws.Range("E14").Value = "  +5.86%  "
$ws.Range("E15").Value = "  +0.59%  "
$ws.Range("D16").Value = "65.481.32"
$ws.Range("E16").Value = "  +6.59%  "
$ws.Range("D17").Value = "3.496.05"
$ws.Range("E17").Value = "  +3.81%  "
$ws.Range("E18").Value = "  +5.56%  "
$ws.Range("D19").Value = "2.990.87"
$ws.Range("E19").Value = "  +3.64%  "
Set-TextValue "D20" "452.65"
$ws.Range("E20").Value = "  +4.78%  "
Set-TextValue "D21" "13.70"
$ws.Range("E21").Value = "  +4.55%  "
$ws.Range("E22").Value = "  +3.62%  "
$ws.Range("E23").Value = "  +7.08%  "
Set-TextValue "D24" "81.32"
$ws.Range("E24").Value = "  +2.58%  "
Set-TextValue "D25" "12.42"
$ws.Range("E25").Value = "  +3.89%  "
Set-TextValue "D26" "2.22"
$ws.Range("E26").Value = "  +10.48%  "
Set-TextValue "D27" "10.64"
$ws.Range("E27").Value = "  +6.23%  "
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("E29").Value = "  +17.22%  "
Set-TextValue "D30" "7.76"
$ws.Range("E30").Value = "  +10.83%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D31" "2.60"
$ws.Range("E31").Value = "  +3.83%  "
$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D32" "0.0000103"
$ws.Range("E32").Value = "  -2.21%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D33" "26.87"
$ws.Range("E33").Value = "  +5.53%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D34" "0.111"
$ws.Range("E34").Value = "  +3.86%  "
Set-TextValue "D35" "1.00"
$ws.Range("E35").Value = "  +0.02%  "
Set-TextValue "D36" "0.985"
$ws.Range("E37").Value = "  +6.90%  "
Set-TextValue "D38" "2.10"
$ws.Range("E38").Value = "  +8.82%  "
Set-TextValue "D39" "45.83"
$ws.Range("E39").Value = "  +17.08%  "
Set-TextValue "D40" "49.18"
$ws.Range("E41").Value = "  +2.35%  "
$ws.Range("E42").Value = "  +5.76%  "
Set-TextValue "D43" "0.298"
$ws.Range("E43").Value = "  +11.48%  "
$ws.Range("E44").Value = "  +2.33%  "
Set-TextValue "D45" "384.03"
$ws.Range("E45").Value = "  +11.76%  "
$ws.Range("D46").Value = "2.766.82"
$ws.Range("E46").Value = "  +2.02%  "
Set-TextValue "D47" "0.0351"
$ws.Range("E47").Value = "  +4.52%  "
Set-TextValue "D48" "135.22"
$ws.Range("E48").Value = "  +1.75%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("E50").Value = "  +2.68%  "
Set-TextValue "D51" "23.16"
$ws.Range("E51").Value = "  +7.29%  "
